$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.610.18'
$ws.Range("E2").Value = '  +0.88%  '

$ws.Range("D3").Value = '3.393.14'
$ws.Range("E3").Value = '  -0.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.63'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.06%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.475'
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.66'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.74%  '

$ws.Range("E10").Value = '  -1.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.388'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.76%  '

$ws.Range("D12").Value = '3.969.72'
$ws.Range("E12").Value = '  -0.11%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.62'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.38%  '

$ws.Range("E14").Value = '  +0.84%  '

$ws.Range("D15").Value = '3.391.08'

$ws.Range("E16").Value = '  -0.75%  '

$ws.Range("D17").Value = '61.543.87'
$ws.Range("E17").Value = '  +0.73%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.18'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.79%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.36%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.98'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.91%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '390.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.95%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.08'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.21%  '

$ws.Range("E23").Value = '  -0.64%  '

$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("E25").Value = '  -4.26%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.191'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.90%  '

$ws.Range("E27").Value = '  -0.15%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.25'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.38%  '

$ws.Range("E29").Value = '  +0.36%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.22%  '

$ws.Range("E31").Value = '  -0.01%  '

$ws.Range("E32").Value = '  -3.86%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.46'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.95'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '167.73'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.42%  '

$ws.Range("D37").Value = '3.425.65'
$ws.Range("E37").Value = '  -0.03%  '

$ws.Range("E38").Value = '  -1.28%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0768'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.32%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.65'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -9.69%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.780'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.08%  '

$ws.Range("E42").Value = '  -0.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.66'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.77%  '

$ws.Range("E44").Value = '  -0.12%  '

$ws.Range("D45").Value = '2.458.61'
$ws.Range("E45").Value = '  -1.26%  '

$ws.Range("E46").Value = '  -2.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.63'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.38%  '

$ws.Range("E48").Value = '  +0.07%  '

$ws.Range("E49").Value = '  -3.86%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.03'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.78%  '

$ws.Range("E51").Value = '  -1.85%  '

